$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 10563
$ws.Range("I74").Value = 12500
$ws.Range("J74").Value = 9917.333000000001
$ws.Range("K74").Value = 12500
$ws.Range("L74").Value = 9917.333000000001
$ws.Range("M74").Value = -11564
$ws.Range("N74").Value = -11789.333
$ws.Range("H77").Value = 10563
$ws.Range("I77").Value = 12500
$ws.Range("J77").Value = 9917.333000000001
$ws.Range("K77").Value = 62500
$ws.Range("L77").Value = 49586.665
$ws.Range("M77").Value = -57820
$ws.Range("N77").Value = -58946.665
$ws.Range("H100").Value = 92106.336
$ws.Range("I100").Value = 10432.429
$ws.Range("K100").Value = 10432.429
$ws.Range("M100").Value = -9891.429
$ws.Range("H107").Value = 3485.8823
$ws.Range("I107").Value = 3585.3333
$ws.Range("K107").Value = 3585.3333
$ws.Range("M107").Value = -1665.3333
$ws.Range("H113").Value = 8221.392
$ws.Range("I113").Value = 10818.833
$ws.Range("K113").Value = 10818.833
$ws.Range("M113").Value = -7564.833000000001
$ws.Range("H115").Value = 867.3333
$ws.Range("I115").Value = 524.75
$ws.Range("J115").Value = 1552.5
$ws.Range("K115").Value = 1574.25
$ws.Range("L115").Value = 4657.5
$ws.Range("M115").Value = -7.25
$ws.Range("N115").Value = -7791.5
$ws.Range("H116").Value = 8728.286
$ws.Range("J116").Value = 10701
$ws.Range("L116").Value = 10701
$ws.Range("N116").Value = -17585
$ws.Range("H137").Value = 3597.7646
$ws.Range("I137").Value = 3091.1365
$ws.Range("J137").Value = 3982.1035
$ws.Range("K137").Value = 9273.4095
$ws.Range("L137").Value = 11946.3105
$ws.Range("M137").Value = -6723.4095
$ws.Range("N137").Value = -17046.3105
$ws.Range("H141").Value = 2067.7334
$ws.Range("I141").Value = 2183.2856
$ws.Range("K141").Value = 6549.8568
$ws.Range("M141").Value = -1369.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1790.1482
$ws.Range("I32").Value = 1833.36
$ws.Range("K32").Value = 1833.36
$ws.Range("M32").Value = -1546.36
$ws.Range("H45").Value = 76927960
$ws.Range("J45").Value = 8232
$ws.Range("L45").Value = 8232
$ws.Range("N45").Value = -8986
$ws.Range("H74").Value = 7755607
$ws.Range("I74").Value = 9526925
$ws.Range("J74").Value = 6089
$ws.Range("K74").Value = 9526925
$ws.Range("L74").Value = 6089
$ws.Range("M74").Value = -9526051
$ws.Range("N74").Value = -7837
$ws.Range("H77").Value = 7755607
$ws.Range("I77").Value = 9526925
$ws.Range("J77").Value = 6089
$ws.Range("K77").Value = 47634625
$ws.Range("L77").Value = 30445
$ws.Range("M77").Value = -47630257
$ws.Range("N77").Value = -39181
$ws.Range("H110").Value = 4401.357
$ws.Range("I110").Value = 2778.7778
$ws.Range("J110").Value = 7322
$ws.Range("K110").Value = 2778.7778
$ws.Range("L110").Value = 7322
$ws.Range("M110").Value = -733.7777999999998
$ws.Range("N110").Value = -11412
$ws.Range("H132").Value = 4516.851
$ws.Range("I132").Value = 3681.425
$ws.Range("J132").Value = 9290.714
$ws.Range("K132").Value = 11044.275
$ws.Range("L132").Value = 27872.142
$ws.Range("M132").Value = -8514.275000000001
$ws.Range("N132").Value = -32932.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2831.0386
$ws.Range("I107").Value = 2428.3333
$ws.Range("K107").Value = 2428.3333
$ws.Range("M107").Value = -508.3332999999998
$ws.Range("H134").Value = 3557.9092
$ws.Range("I134").Value = 1883.5
$ws.Range("J134").Value = 8023
$ws.Range("K134").Value = 5650.5
$ws.Range("L134").Value = 24069
$ws.Range("M134").Value = -3115.5
$ws.Range("N134").Value = -29139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 505.7143
$ws.Range("J5").Value = 647.5
$ws.Range("L5").Value = 647.5
$ws.Range("N5").Value = -871.5
$ws.Range("H16").Value = 2991
$ws.Range("I16").Value = 2209.7778
$ws.Range("J16").Value = 6506.5
$ws.Range("K16").Value = 2209.7778
$ws.Range("L16").Value = 6506.5
$ws.Range("M16").Value = -1922.7778
$ws.Range("N16").Value = -7080.5
$ws.Range("H107").Value = 1424.8889
$ws.Range("I107").Value = 943.3333
$ws.Range("J107").Value = 3832.6667
$ws.Range("K107").Value = 943.3333
$ws.Range("L107").Value = 3832.6667
$ws.Range("M107").Value = 976.6667
$ws.Range("N107").Value = -7672.6667
$ws.Range("H113").Value = 2991
$ws.Range("I113").Value = 2209.7778
$ws.Range("J113").Value = 6506.5
$ws.Range("K113").Value = 2209.7778
$ws.Range("L113").Value = 6506.5
$ws.Range("M113").Value = -39.77779999999984
$ws.Range("N113").Value = -10846.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6344.3477
$ws.Range("I5").Value = 557.1875
$ws.Range("K5").Value = 1671.5625
$ws.Range("M5").Value = -1559.5625
$ws.Range("H113").Value = 1396.3928
$ws.Range("I113").Value = 963.1
$ws.Range("J113").Value = 1637.1111
$ws.Range("K113").Value = 2889.3
$ws.Range("L113").Value = 4911.3333
$ws.Range("M113").Value = -719.3000000000002
$ws.Range("N113").Value = -9251.3333
$ws.Range("H133").Value = 16818
$ws.Range("I133").Value = 16818
$ws.Range("K133").Value = 50454
$ws.Range("M133").Value = -45394
$ws.Range("H135").Value = 6344.3477
$ws.Range("I135").Value = 557.1875
$ws.Range("K135").Value = 5014.6875
$ws.Range("M135").Value = -2479.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1502.2
$ws.Range("I107").Value = 625.8889
$ws.Range("K107").Value = 625.8889
$ws.Range("M107").Value = 1294.1111
$ws.Range("H132").Value = 3029
$ws.Range("I132").Value = 3032.8096
$ws.Range("K132").Value = 9098.4288
$ws.Range("M132").Value = -6568.4288

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4450
$ws.Range("I61").Value = 4450
$ws.Range("K61").Value = 4450
$ws.Range("M61").Value = -4248
$ws.Range("H113").Value = 4450
$ws.Range("I113").Value = 4450
$ws.Range("K113").Value = 4450
$ws.Range("M113").Value = -2280
$ws.Range("H132").Value = 6022.1816
$ws.Range("I132").Value = 3998.1667
$ws.Range("K132").Value = 11994.5001
$ws.Range("M132").Value = -9464.500100000001
$ws.Range("H136").Value = 6246.8945
$ws.Range("I136").Value = 3837.9565
$ws.Range("K136").Value = 11513.8695
$ws.Range("M136").Value = -8963.869499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1109.8889
$ws.Range("I113").Value = 995.0952
$ws.Range("J113").Value = 1270.6
$ws.Range("K113").Value = 2985.2856
$ws.Range("L113").Value = 3811.8
$ws.Range("M113").Value = -815.2856000000002
$ws.Range("N113").Value = -8151.799999999999
$ws.Range("H122").Value = 2960.9312
$ws.Range("I122").Value = 1711.7391
$ws.Range("J122").Value = 7749.5
$ws.Range("K122").Value = 5135.2173
$ws.Range("L122").Value = 23248.5
$ws.Range("M122").Value = -2685.2173
$ws.Range("N122").Value = -28148.5
